$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 12:49:53"
$wsZhCn.Range("H2").Value = "2016-03-17 12:50:16"

# de-de sheet: update Correspond Handoff/Handback Datetime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 12:49:57"
$wsDeDe.Range("H2").Value = "2016-03-17 12:50:31"
